# Apply the commit's changes to the "Test Cases" worksheet:
#  - Clear cell D2 (previously "PASS"), leaving the Results column empty for row 2
#  - Move the active cell / selection from A2 to F20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray "PASS" value from D2
$ws.Range("D2").Value = $null

# Update the current selection shown in the worksheet view
$ws.Range("F20").Select()
